$wb = $excel.ActiveWorkbook

# --- Rename "Sheet2" to "masterTap" ---
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "masterTap"

# --- Populate masterTap with the new list of shared-string values (A1:A8) ---
$ws.Range("A1").Value = "Banks"
$ws.Range("A2").Value = "Users"
$ws.Range("A3").Value = "Branches"
$ws.Range("A4").Value = "Customers"
$ws.Range("A5").Value = "Brokers"
$ws.Range("A6").Value = "Vendors"
$ws.Range("A7").Value = "Enquiry"
$ws.Range("A8").Value = "Projects"

# --- Column A best-fit width (closest reachable value to the recorded 9.6640625) ---
$ws.Columns.Item(1).ColumnWidth = 8.78

# --- Portrait page setup for the new sheet ---
$ws.PageSetup.Orientation = 1

# --- Final UI state: masterTap active with A12 selected ---
$ws.Activate()
$ws.Range("A12").Select()
